$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "G2" = 2.5
    "J2" = 1.11
    "K2" = 6.5
    "G3" = 2.45
    "AA4" = 7.4
    "AD4" = 250
    "AE4" = 10.25
    "AH4" = 18.5
    "G4" = 3.4
    "H4" = 3.7
    "I4" = 1.95
    "J4" = 1.04
    "K4" = 8.75
    "L4" = 1.19
    "M4" = 4.2
    "N4" = 1.57
    "O4" = 2.25
    "P4" = 1.31
    "Q4" = 3.15
    "R4" = 1.52
    "S4" = 2.37
    "T4" = 13.5
    "U4" = 21
    "X4" = 27
    "Y4" = 28
    "Z4" = 8.75
    "G6" = 1.33
    "N6" = 1.87
    "O6" = 1.87
    "G10" = 2.25
    "I10" = 2.82
    "G11" = 1.99
    "H11" = 3.2
    "I11" = 3.5
    "G12" = 2.43
    "H12" = 3.05
    "I12" = 2.75
    "G13" = 1.32
    "H13" = 4.4
    "I13" = 8.800000000000001
    "G14" = 2.75
    "I14" = 2.45
    "G15" = 2.82
    "I15" = 2.57
    "G16" = 2.55
    "I16" = 2.82
    "J16" = 1.1
    "K16" = 7
    "L16" = 1.5
    "M16" = 2.5
    "N16" = 2.6
    "O16" = 1.48
    "AA17" = 7
    "AB17" = 19
    "AH17" = 51
    "G17" = 1.61
    "H17" = 3.75
    "I17" = 5.25
    "N17" = 2.1
    "O17" = 1.7
    "P17" = 1.44
    "Q17" = 2.63
    "Y17" = 34
    "G18" = 2.75
    "I18" = 2.55
    "G19" = 2.65
    "I19" = 2.65
    "G20" = 2.25
    "G21" = 2.82
    "I21" = 2.57
    "G22" = 2.82
    "H22" = 2.82
    "I22" = 2.55
    "G23" = 2.75
    "I23" = 2.2
    "N23" = 1.6
    "O23" = 2.3
    "I24" = 2.05
    "I25" = 2.25
    "K25" = 8
    "I26" = 2.01
    "G27" = 2.65
    "I27" = 2.35
    "G28" = 2.57
    "I28" = 2.57
    "J28" = 1.08
    "K28" = 8
    "G29" = 2.45
    "I29" = 2.65
    "N29" = 1.97
    "O29" = 1.77
    "AB30" = 13
    "AD30" = 201
    "AE30" = 8.5
    "I30" = 2.4
    "K30" = 9.5
    "N30" = 1.97
    "O30" = 1.77
    "R30" = 1.73
    "S30" = 2
    "G31" = 2.05
    "AA32" = 6.3
    "AB32" = 11.5
    "AE32" = 10.75
    "AF32" = 19
    "AI32" = 29
    "AJ32" = 30
    "G32" = 1.72
    "H32" = 3.6
    "I32" = 4
    "L32" = 1.22
    "M32" = 3.85
    "O32" = 1.98
    "R32" = 1.69
    "S32" = 2.04
    "T32" = 6.9
    "U32" = 7.7
    "V32" = 6.9
    "W32" = 11.75
    "X32" = 10.75
    "Z32" = 11.75
    "N33" = 1.77
    "O33" = 1.97
    "G39" = 1.25
    "N40" = 1.94
    "O40" = 1.79
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
